$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (G3/I3): the measure/dimension URIs change from "measure" to "dimension"
$ws.Range("G3").Value = "iaest-dimension:sexo"
$ws.Range("I3").Value = "iaest-dimension:grandes-grupos"

# Row 4 (G4/I4): these columns are now flagged as dimensions ("dim") rather than measures ("medida")
$ws.Range("G4").Value = "dim"
$ws.Range("I4").Value = "dim"

# Row 5 (G5/I5): data type changes from "xsd:string" to "skos:Concept"
$ws.Range("G5").Value = "skos:Concept"
$ws.Range("I5").Value = "skos:Concept"

# New row 6: mapping files for the dimension columns
# (copy the formatting from row 5 so the new cells carry the same style)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$ws.Range("G6").Value = "mapping-sexo.xlsx"
$ws.Range("I6").Value = "mapping-grandes-grupos.xlsx"
